$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing columns C:R to D:S
$ws.Range("C1").EntireColumn.Insert()

# Set the new header for column C
$ws.Range("C1").Value = "TermRun_AGES_year"

# Fill "2023" (as text) for the new column in every data row (rows 2 through 26)
$ws.Range("C2:C26").Formula = '=TEXT(2023,"0")'
$ws.Range("C2:C26").Copy()
$ws.Range("C2:C26").PasteSpecial(-4163)
